# Generate Report for Handoff
# Regenerates the localization-status report for a new handoff pass:
# the source file's GUID-based name changed (b16c76ff... -> d78a5244...),
# its xliff checksums changed, and the handoff timestamps moved forward.

$wb = $excel.ActiveWorkbook

$oldGuidName   = "b16c76ff-8552-40f3-a6f4-5750a8215025.md"
$newGuidName   = "d78a5244-2cef-4e4f-9337-9473f3403082.md"

$oldPathName   = "e2e\b16c76ff-8552-40f3-a6f4-5750a8215025.md"
$newPathName   = "e2e\d78a5244-2cef-4e4f-9337-9473f3403082.md"

$oldHoDate     = "2016-10-25 02:47:10"
$newHoDate     = "2016-10-25 02:47:52"

$oldZhXlf      = "b16c76ff-8552-40f3-a6f4-5750a8215025.115e08a6859a1dffd999b95ea61f60d1ec96eeae.zh-cn.xlf"
$newZhXlf      = "d78a5244-2cef-4e4f-9337-9473f3403082.cc2cc7594cdcb4c433132579631f3ccb66f1d9d7.zh-cn.xlf"

$oldZhDate     = "2016-10-25 02:46:59"
$newZhDate     = "2016-10-25 02:47:40"

$oldDeXlf      = "b16c76ff-8552-40f3-a6f4-5750a8215025.115e08a6859a1dffd999b95ea61f60d1ec96eeae.de-de.xlf"
$newDeXlf      = "d78a5244-2cef-4e4f-9337-9473f3403082.cc2cc7594cdcb4c433132579631f3ccb66f1d9d7.de-de.xlf"

$hyperlinkUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/337f621101861423a9d57326fc868cb90d35cf17/e2e/b16c76ff-8552-40f3-a6f4-5750a8215025.md"

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newHoDate

# Replace the B2 hyperlink so its display text matches the new path, while
# keeping the same target URL.
$linkRange = $wsOverview.Range("B2")
$linkRange.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($linkRange, $hyperlinkUrl, "", "", $newPathName)

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuidName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

$linkRangeZh = $wsZhCn.Range("A2")
$linkRangeZh.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($linkRangeZh, $hyperlinkUrl, "", "", $newGuidName)

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuidName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

$linkRangeDe = $wsDeDe.Range("A2")
$linkRangeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($linkRangeDe, $hyperlinkUrl, "", "", $newGuidName)
